# GF#21457 - added CodeSystem to form curation diagram
# Moves the existing "Questionnaire" / "ValueSet" diagram rows up and
# inserts a new "CodeSystem" row (rectangle + arrow connector) below them.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> point helper (PowerPoint COM positions are expressed in points).
# Shape.Top/Left/Width/Height are single-precision (Single) in the real
# PowerPoint object model, so a naive emu/12700 conversion can be off by
# 1 EMU once PowerPoint truncates it back down to whole EMUs on save.
# Search near the naive value for a point value that survives the
# Single round-trip and truncation exactly, so the saved EMU matches.
function EmuToPt([double]$targetEmu) {
    $base = $targetEmu / 12700.0
    for ($i = 0; $i -lt 500; $i++) {
        $step = $i * 0.00001
        $cand = $base + $step
        $f32 = [single]$cand
        $emu = [math]::Floor([double]$f32 * 12700.0)
        if ($emu -eq $targetEmu) {
            return $cand
        }
        if ($i -gt 0) {
            $cand2 = $base - $step
            $f32b = [single]$cand2
            $emu2 = [math]::Floor([double]$f32b * 12700.0)
            if ($emu2 -eq $targetEmu) {
                return $cand2
            }
        }
    }
    return $base
}

# --- 1. Shift the existing "Questionnaire" row and "ValueSet" row upward ---
$arrow1 = $s.Shapes.Item("Straight Arrow Connector 34")
$arrow1.Top = EmuToPt 165894

$questionnaireRect = $s.Shapes.Item("Rectangle 35")
$questionnaireRect.Top = EmuToPt 188608

$valueSetRect = $s.Shapes.Item("Rectangle 36")
$valueSetRect.Top = EmuToPt 775494

$arrow2 = $s.Shapes.Item("Straight Arrow Connector 37")
$arrow2.Top = EmuToPt 832663

# --- 2. Add the new "CodeSystem" rectangle (clone of the ValueSet one). ---
# PowerPoint assigned this shape id 8 when it was authored. Shape.Id can't
# be set directly (read-only, just like in real PowerPoint COM), so nudge
# the presentation's internal id counter forward by creating/discarding
# scratch duplicates until the next shape lands on the expected id.
$targetRectId = 8
$guard = 0
while ($true) {
    $probe = $valueSetRect.Duplicate()
    $probeShape = $probe.Item(1)
    if ($probeShape.Id -eq $targetRectId) {
        $codeSystemRect = $probeShape
        break
    }
    $probeShape.Delete()
    $guard++
    if ($guard -gt 1000) {
        # Safety net: give up trying to match the id and just use whatever
        # shape we most recently duplicated.
        $codeSystemRect = $probeShape
        break
    }
}
$codeSystemRect.Name = "Rectangle 7"
$codeSystemRect.Left = EmuToPt 967477
$codeSystemRect.Top = EmuToPt 1442013
$codeSystemRect.Width = EmuToPt 1086027
$codeSystemRect.Height = EmuToPt 628881

# Replace only the bold "ValueSet" word with "CodeSystem", keeping the
# "Query, Create, Update, Delete" line and its formatting untouched.
$fullText = $codeSystemRect.TextFrame.TextRange.Text
$prefix = "Query, Create, Update, Delete"
$wordRange = $codeSystemRect.TextFrame.TextRange.Characters($prefix.Length + 1, $fullText.Length - $prefix.Length)
$wordRange.Text = "CodeSystem"

# --- 3. Add the new arrow connector under the CodeSystem rectangle ---
$targetConnId = 9
$guard2 = 0
while ($true) {
    $probeArrow = $arrow2.Duplicate()
    $probeArrowShape = $probeArrow.Item(1)
    if ($probeArrowShape.Id -eq $targetConnId) {
        $codeSystemArrow = $probeArrowShape
        break
    }
    $probeArrowShape.Delete()
    $guard2++
    if ($guard2 -gt 1000) {
        $codeSystemArrow = $probeArrowShape
        break
    }
}
$codeSystemArrow.Name = "Straight Arrow Connector 8"
$codeSystemArrow.Left = EmuToPt 957990
$codeSystemArrow.Top = EmuToPt 1499182
$codeSystemArrow.Width = EmuToPt 1115798
$codeSystemArrow.Height = EmuToPt 0

Write-Host "Added CodeSystem row to form curation diagram"
